# Automatische test-sync: 2025-08-05 16:52:50
# Append the new "Kun jij dit even regelen?" test-mail log entry (row 9) to the
# "Logs" sheet, extend the conditional formatting ranges to include it, and
# bump the "Planning / Afspraak" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- New row 9 in the Logs sheet -------------------------------------------------
$wsLogs.Range("A9").Value = "Kun jij dit even regelen?"
$wsLogs.Range("B9").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C9").Value = "Testmail #1: Kun jij dit even regelen?"
$wsLogs.Range("D9").Value = "Planning / Afspraak"
$wsLogs.Range("E9").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$wsLogs.Range("F9").Value = "2025-08-05 16:52:35"
$wsLogs.Range("G9").Value = "Ja"
$wsLogs.Range("H9").Value = "Ja"
$wsLogs.Range("I9").Value = "Nee"
$wsLogs.Range("J9").Value = "Nee"

# --- Extend the conditional formatting ranges (D/G/H/I/J) down to row 9 ---------
$wsLogs.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D9"))
$wsLogs.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G9"))
$wsLogs.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H9"))
$wsLogs.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I9"))
$wsLogs.Range("J2:J8").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J9"))

# --- Update the Dashboard count for "Planning / Afspraak" -----------------------
$wsDash.Range("B3").Value = 3
